$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Review")

function Copy-CellFormat($srcAddr, $dstAddr) {
    $ws.Range($srcAddr).Copy() | Out-Null
    $ws.Range($dstAddr).PasteSpecial(-4122) | Out-Null  # xlPasteFormats
}

# --- Row 4: new review entry (Tutorial_03) ---
$ws.Range("F4").Value = "Tutorial_03`nindex.php"
$ws.Range("I4").Value = @"
1) Use an indent of 4 spaces for *.php
2) When enter future date,instead of showing your age is 0,
   use condition to show message for Invalid DOB.
3) Use proper name 
"@

# --- Row 2 / Row 3: "change 4 sapce" -> "change 4 space", and response dates shift later the same day ---
$ws.Range("P2").Value = "change 4 space"
$ws.Range("V2").Value = 44546.583333333336

$ws.Range("P3").Value = "change 4 space"
$ws.Range("V3").Value = 44546.583333333336

$ws.Range("O4").Value = "PyaePyaeHan"
$ws.Range("P4").Value = @"
1) change 4 space
2)  show message  DOB.
3) Use proper name 
"@
$ws.Range("V4").Value = 44547.833333333336
$ws.Range("W4").Value = "HeinHtetSan"

Copy-CellFormat "F2" "F4"
Copy-CellFormat "G2" "G4"
Copy-CellFormat "H2" "H4"
Copy-CellFormat "F2" "P4"
Copy-CellFormat "G2" "Q4"
Copy-CellFormat "G2" "R4"
Copy-CellFormat "G2" "S4"
Copy-CellFormat "G2" "T4"
Copy-CellFormat "H2" "U4"
Copy-CellFormat "V2" "V4"

# --- Row 5: new review entry (Tutorial_04) ---
$ws.Range("F5").Value = "Tutorial_04`nindex.php"

$ws.Range("O5").Value = "PyaePyaeHan"
$ws.Range("P5").Value = @"
1) change 4 space
2)  set connection Login
3) connectin not accept null data
"@
$ws.Range("V5").Value = 44547.833333333336
$ws.Range("W5").Value = "HeinHtetSan"

$ws.Range("I5").Value = @"
1) Use an indent of 4 spaces for *.php
2) Please set default username and password to show Invalid Login
3) Need form validation (not form accepting null data)
"@

Copy-CellFormat "F2" "P5"
Copy-CellFormat "G2" "Q5"
Copy-CellFormat "G2" "R5"
Copy-CellFormat "G2" "S5"
Copy-CellFormat "G2" "T5"
Copy-CellFormat "H2" "U5"
Copy-CellFormat "V2" "V5"

$excel.CutCopyMode = 0

# --- Selection / active view state ---
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollColumn = 6
$win.ScrollRow = 4
$ws.Range("I5:N5").Select()

Write-Output "Applied review updates for Tutorial_03 and Tutorial_04"
